$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 37768.9
$ws.Range("J3").Value = 37768.9
$ws.Range("L3").Value = 37768.9
$ws.Range("N3").Value = -37996.9
$ws.Range("H93").Value = 49601
$ws.Range("J93").Value = 49601
$ws.Range("L93").Value = 49601
$ws.Range("N93").Value = -54593
$ws.Range("H95").Value = 31638
$ws.Range("J95").Value = 31638
$ws.Range("L95").Value = 31638
$ws.Range("N95").Value = -37130
$ws.Range("H102").Value = 37768.9
$ws.Range("J102").Value = 37768.9
$ws.Range("L102").Value = 37768.9
$ws.Range("N102").Value = -44258.9
$ws.Range("H105").Value = 48896
$ws.Range("J105").Value = 48896
$ws.Range("L105").Value = 48896
$ws.Range("N105").Value = -55884

# Sheet 2: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 31870
$ws.Range("J24").Value = 31870
$ws.Range("L24").Value = 31870
$ws.Range("N24").Value = -32618
$ws.Range("H32").Value = 34924.758
$ws.Range("I32").Value = 33633.402
$ws.Range("J32").Value = 68500
$ws.Range("K32").Value = 33633.402
$ws.Range("L32").Value = 68500
$ws.Range("M32").Value = -33346.402
$ws.Range("N32").Value = -69074
$ws.Range("H70").Value = 20000
$ws.Range("J70").Value = 20000
$ws.Range("L70").Value = 20000
$ws.Range("N70").Value = -20540
$ws.Range("H73").Value = 20000
$ws.Range("J73").Value = 20000
$ws.Range("L73").Value = 20000
$ws.Range("N73").Value = -21872
$ws.Range("H94").Value = 33300
$ws.Range("J94").Value = 33300
$ws.Range("L94").Value = 33300
$ws.Range("N94").Value = -35102
$ws.Range("H96").Value = 32897
$ws.Range("J96").Value = 32897
$ws.Range("L96").Value = 32897
$ws.Range("N96").Value = -38389
$ws.Range("H100").Value = 31870
$ws.Range("J100").Value = 31870
$ws.Range("L100").Value = 31870
$ws.Range("N100").Value = -34034
$ws.Range("H101").Value = 45760.4
$ws.Range("J101").Value = 45760.4
$ws.Range("L101").Value = 45760.4
$ws.Range("N101").Value = -52250.4
$ws.Range("H103").Value = 37794.668
$ws.Range("J103").Value = 37794.668
$ws.Range("L103").Value = 37794.668
$ws.Range("N103").Value = -40138.668
$ws.Range("H105").Value = 48868
$ws.Range("J105").Value = 48868
$ws.Range("L105").Value = 48868
$ws.Range("N105").Value = -55856
$ws.Range("H106").Value = 46445.332
$ws.Range("J106").Value = 46445.332
$ws.Range("L106").Value = 46445.332
$ws.Range("N106").Value = -48969.332

# Sheet 3: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 20000
$ws.Range("J21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20472
$ws.Range("H92").Value = 33467
$ws.Range("J92").Value = 33467
$ws.Range("L92").Value = 33467
$ws.Range("N92").Value = -38459
$ws.Range("H95").Value = 43616
$ws.Range("J95").Value = 43616
$ws.Range("L95").Value = 43616
$ws.Range("N95").Value = -49108
$ws.Range("H103").Value = 199497.5
$ws.Range("J103").Value = 199497.5
$ws.Range("L103").Value = 199497.5
$ws.Range("N103").Value = -201841.5
$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 30000
$ws.Range("N106").Value = -32524

# Sheet 4: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 47996
$ws.Range("J43").Value = 47996
$ws.Range("L43").Value = 47996
$ws.Range("N43").Value = -48364
$ws.Range("H92").Value = 46601
$ws.Range("J92").Value = 46601
$ws.Range("L92").Value = 46601
$ws.Range("N92").Value = -51593
$ws.Range("H96").Value = 79498
$ws.Range("J96").Value = 79498
$ws.Range("L96").Value = 79498
$ws.Range("N96").Value = -84990
$ws.Range("H101").Value = 47996
$ws.Range("J101").Value = 47996
$ws.Range("L101").Value = 47996
$ws.Range("N101").Value = -54486
$ws.Range("H106").Value = 19992
$ws.Range("J106").Value = 19992
$ws.Range("L106").Value = 19992
$ws.Range("N106").Value = -22516
$ws.Range("H131").Value = 38326
$ws.Range("J131").Value = 38326
$ws.Range("L131").Value = 38326
$ws.Range("N131").Value = -48406

# Sheet 5: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 76282.32000000001
$ws.Range("I131").Value = 14722.714
$ws.Range("J131").Value = 96802.19
$ws.Range("K131").Value = 44168.142
$ws.Range("L131").Value = 290406.57
$ws.Range("M131").Value = -39128.142
$ws.Range("N131").Value = -300486.57

# Sheet 6: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 33360134
$ws.Range("J98").Value = 33360134
$ws.Range("L98").Value = 33360134
$ws.Range("N98").Value = -33366124
$ws.Range("H100").Value = 35320
$ws.Range("J100").Value = 35320
$ws.Range("L100").Value = 35320
$ws.Range("N100").Value = -37484
$ws.Range("H104").Value = 38998
$ws.Range("J104").Value = 38998
$ws.Range("L104").Value = 38998
$ws.Range("N104").Value = -45986
$ws.Range("H105").Value = 46996
$ws.Range("J105").Value = 46996
$ws.Range("L105").Value = 46996
$ws.Range("N105").Value = -53984
$ws.Range("H125").Value = 44326
$ws.Range("J125").Value = 44326
$ws.Range("L125").Value = 44326
$ws.Range("N125").Value = -49246

# Sheet 7: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 44389
$ws.Range("J92").Value = 44389
$ws.Range("L92").Value = 44389
$ws.Range("N92").Value = -49381
$ws.Range("H94").Value = 43211.223
$ws.Range("J94").Value = 43211.223
$ws.Range("L94").Value = 43211.223
$ws.Range("N94").Value = -44563.223
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H98").Value = 39000
$ws.Range("J98").Value = 39000
$ws.Range("L98").Value = 39000
$ws.Range("N98").Value = -44990
$ws.Range("H103").Value = 48590
$ws.Range("J103").Value = 48590
$ws.Range("L103").Value = 48590
$ws.Range("N103").Value = -50934
$ws.Range("H105").Value = 49615
$ws.Range("J105").Value = 49615
$ws.Range("L105").Value = 49615
$ws.Range("N105").Value = -56603
$ws.Range("H106").Value = 46000
$ws.Range("J106").Value = 46000
$ws.Range("L106").Value = 46000
$ws.Range("N106").Value = -48524
$ws.Range("H117").Value = 45056
$ws.Range("J117").Value = 45056
$ws.Range("L117").Value = 45056
$ws.Range("N117").Value = -54234
$ws.Range("H123").Value = 42177
$ws.Range("J123").Value = 42177
$ws.Range("L123").Value = 42177
$ws.Range("N123").Value = -51977
$ws.Range("H129").Value = 37929
$ws.Range("J129").Value = 37929
$ws.Range("L129").Value = 37929
$ws.Range("N129").Value = -47929

# Sheet 8: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 29832.5
$ws.Range("J92").Value = 29832.5
$ws.Range("L92").Value = 29832.5
$ws.Range("N92").Value = -34824.5
$ws.Range("H98").Value = 43780.5
$ws.Range("J98").Value = 43780.5
$ws.Range("L98").Value = 43780.5
$ws.Range("N98").Value = -49770.5
$ws.Range("H103").Value = 42148
$ws.Range("J103").Value = 42148
$ws.Range("L103").Value = 42148
$ws.Range("N103").Value = -44492
$ws.Range("H104").Value = 40480.668
$ws.Range("J104").Value = 40480.668
$ws.Range("L104").Value = 40480.668
$ws.Range("N104").Value = -47468.668
$ws.Range("H109").Value = 32406.5
$ws.Range("J109").Value = 32406.5
$ws.Range("L109").Value = 32406.5
$ws.Range("N109").Value = -35180.5
$ws.Range("H127").Value = 33143
$ws.Range("J127").Value = 33143
$ws.Range("L127").Value = 33143
$ws.Range("N127").Value = -43063
$ws.Range("H129").Value = 30516.455
$ws.Range("J129").Value = 30516.455
$ws.Range("L129").Value = 30516.455
$ws.Range("N129").Value = -40516.455
